$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells in column D sometimes hold numeric-looking text (e.g. "247.76")
# that Excel would otherwise auto-convert to a real number. Pre-formatting the
# cell as Text ("@") before the assignment keeps it a literal string, matching
# the source data which always stores prices as text.

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "36.663.06"
$ws.Cells.Item(2, 5).Value = "  +0.78%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "2.008.17"
$ws.Cells.Item(3, 5).Value = "  -0.25%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  +0.06%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "247.76"
$ws.Cells.Item(5, 5).Value = "  -1.78%  "

# Row 6 - XRP
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.632"
$ws.Cells.Item(6, 5).Value = "  -1.47%  "

# Row 7 - Solana
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "62.35"
$ws.Cells.Item(7, 5).Value = "  +0.13%  "

# Row 9 - Cardano
$ws.Cells.Item(9, 5).Value = "  +3.63%  "

# Row 10 - OKB
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "57.30"
$ws.Cells.Item(10, 5).Value = "  -1.59%  "

# Row 11 - Dogecoin
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0781"
$ws.Cells.Item(11, 5).Value = "  +5.04%  "

# Row 12 - TRON
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.103"
$ws.Cells.Item(12, 5).Value = "  -0.29%  "

# Row 13 - Polygon
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.887"
$ws.Cells.Item(13, 5).Value = "  -1.69%  "

# Row 14 - Avalanche
$ws.Cells.Item(14, 5).Value = "  +8.71%  "

# Row 15 - Chainlink
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "14.21"
$ws.Cells.Item(15, 5).Value = "  -4.69%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Cells.Item(16, 4).Value = "2.303.60"
$ws.Cells.Item(16, 5).Value = "  -0.19%  "

# Row 17 - Polkadot
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "5.52"
$ws.Cells.Item(17, 5).Value = "  +1.07%  "

# Row 18 - WrappedEther
$ws.Cells.Item(18, 4).Value = "2.009.04"
$ws.Cells.Item(18, 5).Value = "  -0.43%  "

# Row 19 - WrappedBTC
$ws.Cells.Item(19, 4).Value = "36.579.17"
$ws.Cells.Item(19, 5).Value = "  +0.64%  "

# Row 20 - Litecoin
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "71.94"
$ws.Cells.Item(20, 5).Value = "  -0.14%  "

# Row 21 - ShibaInu
$ws.Cells.Item(21, 5).Value = "  +0.91%  "

# Row 22 - Uniswap
$ws.Cells.Item(22, 5).Value = "  +0.06%  "

# Row 23 - BitcoinCash
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "237.99"
$ws.Cells.Item(23, 5).Value = "  +1.47%  "

# Row 24 - Dai
$ws.Cells.Item(24, 5).Value = "  -0.02%  "

# Row 25 - PancakeSwap
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.51"
$ws.Cells.Item(25, 5).Value = "  -7.02%  "

# Row 26 - Toncoin
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.32"
$ws.Cells.Item(26, 5).Value = "  -0.08%  "

# Row 27 - Cosmos
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "9.86"
$ws.Cells.Item(27, 5).Value = "  +2.42%  "

# Row 28 - Kaspa
$ws.Cells.Item(28, 5).Value = "  +27.03%  "

# Row 29 - Monero
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "160.11"
$ws.Cells.Item(29, 5).Value = "  -2.02%  "

# Row 30 - EthereumClassic
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "20.11"
$ws.Cells.Item(30, 5).Value = "  +2.41%  "

# Row 31 - Stellar
$ws.Cells.Item(31, 5).Value = "  +0.55%  "

# Row 32 - ImmutableX
$ws.Cells.Item(32, 5).Value = "  -0.41%  "

# Row 33 - Filecoin
$ws.Cells.Item(33, 5).Value = "  -2.55%  "

# Row 34 - Hedera
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0623"
$ws.Cells.Item(34, 5).Value = "  +2.41%  "

# Row 35 - InternetComputer(DFINITY)
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.47"
$ws.Cells.Item(35, 5).Value = "  -2.06%  "

# Row 36 - THORChain
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "6.54"
$ws.Cells.Item(36, 5).Value = "  +10.08%  "

# Row 37 - LidoDAOToken
$ws.Cells.Item(37, 5).Value = "  -4.04%  "

# Row 38 - BinanceUSD
$ws.Cells.Item(38, 5).Value = "  +0.17%  "

# Row 39 - WEMIXToken
$ws.Cells.Item(39, 5).Value = "  +0.79%  "

# Row 40 - RenderToken
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.20"
$ws.Cells.Item(40, 5).Value = "  +20.55%  "

# Row 41 - TrustWalletToken
$ws.Cells.Item(41, 5).Value = "  +2.97%  "

# Row 42 - Cronos
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.100"
$ws.Cells.Item(42, 5).Value = "  -5.06%  "

# Row 43 - HuobiToken
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.93"
$ws.Cells.Item(43, 5).Value = "  +0.38%  "

# Rows 44 & 45 swap place: VeChain moves to row 44, ARBITRUM moves to row 45
# Row 44 becomes VeChain
$ws.Cells.Item(44, 2).Value = "VeChain"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0215"
$ws.Cells.Item(44, 5).Value = "  -0.80%  "

# Row 45 becomes ARBITRUM
$ws.Cells.Item(45, 2).Value = "ARBITRUM"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.12"
$ws.Cells.Item(45, 5).Value = "  -0.97%  "

# Row 46 - InjectiveProtocol
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "16.71"
$ws.Cells.Item(46, 5).Value = "  -2.21%  "

# Row 47 - Aave
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "93.42"
$ws.Cells.Item(47, 5).Value = "  -1.74%  "

# Row 48 - FraxShare
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "7.65"
$ws.Cells.Item(48, 5).Value = "  -3.37%  "

# Row 49 - Maker
$ws.Cells.Item(49, 4).Value = "1.359.21"
$ws.Cells.Item(49, 5).Value = "  -6.46%  "

# Row 50 - MXToken
$ws.Cells.Item(50, 5).Value = "  -1.47%  "

# Row 51 - RocketPoolETH
$ws.Cells.Item(51, 4).Value = "2.194.99"
$ws.Cells.Item(51, 5).Value = "  -0.06%  "
